# Update FFXIV Leve market-profit columns (H-N) per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$alc = $wb.Worksheets.Item("ALC")
$arm = $wb.Worksheets.Item("ARM")
$bsm = $wb.Worksheets.Item("BSM")
$crp = $wb.Worksheets.Item("CRP")
$cul = $wb.Worksheets.Item("CUL")
$gsm = $wb.Worksheets.Item("GSM")
$ltw = $wb.Worksheets.Item("LTW")
$wvr = $wb.Worksheets.Item("WVR")

# ALC row 132
$alc.Cells.Item(132, 8).Value = 1856889.8
$alc.Cells.Item(132, 9).Value = 2269242.2
$alc.Cells.Item(132, 10).Value = 1303.7858
$alc.Cells.Item(132, 11).Value = 6807726.600000001
$alc.Cells.Item(132, 12).Value = 3911.3574
$alc.Cells.Item(132, 13).Value = -6805196.600000001
$alc.Cells.Item(132, 14).Value = -8971.357400000001

# ALC row 141
$alc.Cells.Item(141, 8).Value = 2397.651
$alc.Cells.Item(141, 10).Value = 3634.742
$alc.Cells.Item(141, 12).Value = 10904.226
$alc.Cells.Item(141, 14).Value = -21264.226

# ARM row 2
$arm.Cells.Item(2, 8).Value = 1327.3334
$arm.Cells.Item(2, 9).Value = 1008.1
$arm.Cells.Item(2, 11).Value = 1008.1
$arm.Cells.Item(2, 13).Value = -895.1

# ARM row 61
$arm.Cells.Item(61, 8).Value = 1370.125
$arm.Cells.Item(61, 9).Value = 863.9783
$arm.Cells.Item(61, 10).Value = 3698.4
$arm.Cells.Item(61, 11).Value = 863.9783
$arm.Cells.Item(61, 12).Value = 3698.4
$arm.Cells.Item(61, 13).Value = -651.9783
$arm.Cells.Item(61, 14).Value = -4122.4

# ARM row 116
$arm.Cells.Item(116, 8).Value = 1327.3334
$arm.Cells.Item(116, 9).Value = 1008.1
$arm.Cells.Item(116, 11).Value = 1008.1
$arm.Cells.Item(116, 13).Value = 1285.9

# ARM row 132
$arm.Cells.Item(132, 8).Value = 3348.818
$arm.Cells.Item(132, 9).Value = 3423.0435
$arm.Cells.Item(132, 10).Value = 2969.4443
$arm.Cells.Item(132, 11).Value = 10269.1305
$arm.Cells.Item(132, 12).Value = 8908.332900000001
$arm.Cells.Item(132, 13).Value = -7739.130500000001
$arm.Cells.Item(132, 14).Value = -13968.3329

# ARM row 136
$arm.Cells.Item(136, 8).Value = 1370.125
$arm.Cells.Item(136, 9).Value = 863.9783
$arm.Cells.Item(136, 10).Value = 3698.4
$arm.Cells.Item(136, 11).Value = 2591.9349
$arm.Cells.Item(136, 12).Value = 11095.2
$arm.Cells.Item(136, 13).Value = -41.9349000000002
$arm.Cells.Item(136, 14).Value = -16195.2

# BSM row 3
$bsm.Cells.Item(3, 8).Value = 1327.3334
$bsm.Cells.Item(3, 9).Value = 1008.1
$bsm.Cells.Item(3, 11).Value = 1008.1
$bsm.Cells.Item(3, 13).Value = -894.1

# BSM row 94
$bsm.Cells.Item(94, 8).Value = 685.825
$bsm.Cells.Item(94, 9).Value = 483.55173
$bsm.Cells.Item(94, 10).Value = 1219.091
$bsm.Cells.Item(94, 11).Value = 483.55173
$bsm.Cells.Item(94, 12).Value = 1219.091
$bsm.Cells.Item(94, 13).Value = -32.55173000000002
$bsm.Cells.Item(94, 14).Value = -2121.091

# BSM row 134
$bsm.Cells.Item(134, 8).Value = 19225.316
$bsm.Cells.Item(134, 9).Value = 24600.861
$bsm.Cells.Item(134, 10).Value = 2714.7144
$bsm.Cells.Item(134, 11).Value = 73802.583
$bsm.Cells.Item(134, 12).Value = 8144.1432
$bsm.Cells.Item(134, 13).Value = -71267.583
$bsm.Cells.Item(134, 14).Value = -13214.1432

# CRP row 31
$crp.Cells.Item(31, 8).Value = 3146765.5
$crp.Cells.Item(31, 9).Value = 1385.4667
$crp.Cells.Item(31, 10).Value = 20839528
$crp.Cells.Item(31, 11).Value = 1385.4667
$crp.Cells.Item(31, 12).Value = 20839528
$crp.Cells.Item(31, 13).Value = -1090.4667
$crp.Cells.Item(31, 14).Value = -20840118

# CRP row 34
$crp.Cells.Item(34, 8).Value = 3146765.5
$crp.Cells.Item(34, 9).Value = 1385.4667
$crp.Cells.Item(34, 10).Value = 20839528
$crp.Cells.Item(34, 11).Value = 1385.4667
$crp.Cells.Item(34, 12).Value = 20839528
$crp.Cells.Item(34, 13).Value = -1183.4667
$crp.Cells.Item(34, 14).Value = -20839932

# CRP row 58
$crp.Cells.Item(58, 8).Value = 7247145
$crp.Cells.Item(58, 9).Value = 763.3
$crp.Cells.Item(58, 10).Value = 55556356
$crp.Cells.Item(58, 11).Value = 763.3
$crp.Cells.Item(58, 12).Value = 55556356
$crp.Cells.Item(58, 13).Value = -560.3
$crp.Cells.Item(58, 14).Value = -55556762

# CRP row 62
$crp.Cells.Item(62, 8).Value = 19610720
$crp.Cells.Item(62, 9).Value = 2498.3333
$crp.Cells.Item(62, 10).Value = 41669970
$crp.Cells.Item(62, 11).Value = 2498.3333
$crp.Cells.Item(62, 12).Value = 41669970
$crp.Cells.Item(62, 13).Value = -1874.3333
$crp.Cells.Item(62, 14).Value = -41671218

# CRP row 65
$crp.Cells.Item(65, 8).Value = 19610720
$crp.Cells.Item(65, 9).Value = 2498.3333
$crp.Cells.Item(65, 10).Value = 41669970
$crp.Cells.Item(65, 11).Value = 12491.6665
$crp.Cells.Item(65, 12).Value = 208349850
$crp.Cells.Item(65, 13).Value = -9371.666499999999
$crp.Cells.Item(65, 14).Value = -208356090

# CRP row 132
$crp.Cells.Item(132, 8).Value = 1587.4058
$crp.Cells.Item(132, 9).Value = 1413.3158
$crp.Cells.Item(132, 10).Value = 2414.3333
$crp.Cells.Item(132, 11).Value = 4239.9474
$crp.Cells.Item(132, 12).Value = 7242.999899999999
$crp.Cells.Item(132, 13).Value = -1709.9474
$crp.Cells.Item(132, 14).Value = -12302.9999

# CRP row 134
$crp.Cells.Item(134, 8).Value = 840.5
$crp.Cells.Item(134, 9).Value = 790
$crp.Cells.Item(134, 10).Value = 1168.75
$crp.Cells.Item(134, 11).Value = 2370
$crp.Cells.Item(134, 12).Value = 3506.25
$crp.Cells.Item(134, 13).Value = 165
$crp.Cells.Item(134, 14).Value = -8576.25

# CRP row 136
$crp.Cells.Item(136, 8).Value = 7247145
$crp.Cells.Item(136, 9).Value = 763.3
$crp.Cells.Item(136, 10).Value = 55556356
$crp.Cells.Item(136, 11).Value = 2289.9
$crp.Cells.Item(136, 12).Value = 166669068
$crp.Cells.Item(136, 13).Value = 260.1000000000004
$crp.Cells.Item(136, 14).Value = -166674168

# CUL row 107
$cul.Cells.Item(107, 8).Value = 409.46667
$cul.Cells.Item(107, 10).Value = 419.1111
$cul.Cells.Item(107, 12).Value = 1257.3333
$cul.Cells.Item(107, 14).Value = -5097.3333

# CUL row 140
$cul.Cells.Item(140, 8).Value = 1939.2858
$cul.Cells.Item(140, 9).Value = 1742.3077
$cul.Cells.Item(140, 11).Value = 5226.9231
$cul.Cells.Item(140, 13).Value = -46.92309999999998

# GSM row 70
$gsm.Cells.Item(70, 8).Value = 9448413
$gsm.Cells.Item(70, 9).Value = 11864446
$gsm.Cells.Item(70, 10).Value = 3918.182
$gsm.Cells.Item(70, 11).Value = 11864446
$gsm.Cells.Item(70, 12).Value = 3918.182
$gsm.Cells.Item(70, 13).Value = -11864176
$gsm.Cells.Item(70, 14).Value = -4458.182

# GSM row 73
$gsm.Cells.Item(73, 8).Value = 9448413
$gsm.Cells.Item(73, 9).Value = 11864446
$gsm.Cells.Item(73, 10).Value = 3918.182
$gsm.Cells.Item(73, 11).Value = 11864446
$gsm.Cells.Item(73, 12).Value = 3918.182
$gsm.Cells.Item(73, 13).Value = -11863510
$gsm.Cells.Item(73, 14).Value = -5790.182

# LTW row 40
$ltw.Cells.Item(40, 8).Value = 872.9
$ltw.Cells.Item(40, 9).Value = 877.3889
$ltw.Cells.Item(40, 10).Value = 866.1667
$ltw.Cells.Item(40, 11).Value = 877.3889
$ltw.Cells.Item(40, 12).Value = 866.1667
$ltw.Cells.Item(40, 13).Value = -741.3889
$ltw.Cells.Item(40, 14).Value = -1138.1667

# LTW row 132
$ltw.Cells.Item(132, 8).Value = 5917.915
$ltw.Cells.Item(132, 9).Value = 7486.5625
$ltw.Cells.Item(132, 11).Value = 22459.6875
$ltw.Cells.Item(132, 13).Value = -19929.6875

# WVR row 122
$wvr.Cells.Item(122, 8).Value = 30837.383
$wvr.Cells.Item(122, 9).Value = 39779.69
$wvr.Cells.Item(122, 10).Value = 1774.875
$wvr.Cells.Item(122, 11).Value = 119339.07
$wvr.Cells.Item(122, 12).Value = 5324.625
$wvr.Cells.Item(122, 13).Value = -116889.07
$wvr.Cells.Item(122, 14).Value = -10224.625

# WVR row 132
$wvr.Cells.Item(132, 8).Value = 989.19446
$wvr.Cells.Item(132, 9).Value = 844.0164
$wvr.Cells.Item(132, 10).Value = 1794.2727
$wvr.Cells.Item(132, 11).Value = 2532.0492
$wvr.Cells.Item(132, 12).Value = 5382.8181
$wvr.Cells.Item(132, 13).Value = -2.049199999999928
$wvr.Cells.Item(132, 14).Value = -10442.8181

# WVR row 136
$wvr.Cells.Item(136, 8).Value = 2358.16
$wvr.Cells.Item(136, 9).Value = 2603.6167
$wvr.Cells.Item(136, 10).Value = 1376.3334
$wvr.Cells.Item(136, 11).Value = 7810.8501
$wvr.Cells.Item(136, 12).Value = 4129.0002
$wvr.Cells.Item(136, 13).Value = -5260.8501
$wvr.Cells.Item(136, 14).Value = -9229.0002
